$d = $word.ActiveDocument

$replacements = @(
    @{old="55×20="; new="83×13="},
    @{old="64×80="; new="83×25="},
    @{old="76×92="; new="28×41="},
    @{old="83×37="; new="15×11="},
    @{old="28×49="; new="59×23="},
    @{old="44×18="; new="62×56="},
    @{old="63×92="; new="52×54="},
    @{old="18×57="; new="81×97="},
    @{old="15×71="; new="62×44="},
    @{old="37×49="; new="85×99="},
    @{old="25×26="; new="56×57="},
    @{old="64×31="; new="57×30="},
    @{old="89×59="; new="20×61="},
    @{old="76×58="; new="69×70="},
    @{old="22×31="; new="13×80="},
    @{old="26×56="; new="39×82="},
    @{old="46×35="; new="25×78="},
    @{old="56×21="; new="68×73="},
    @{old="24×54="; new="44×11="},
    @{old="39×43="; new="21×31="},
    @{old="71×24="; new="82×70="},
    @{old="48×88="; new="86×79="},
    @{old="45×51="; new="68×73="},
    @{old="59×85="; new="35×68="},
    @{old="51×34="; new="18×16="}
)

foreach ($r in $replacements) {
    $find = $d.Content.Find
    $find.ClearFormatting()
    $find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
